$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column N into the new column O for rows 4-14,
# then set the appropriate values.
$ws.Range("N4:N14").Copy()
$ws.Range("O4:O14").PasteSpecial(-4122)

$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 2
$ws.Range("O6").Value = "-"
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = "-"
$ws.Range("O9").Value = "-"
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = "-"
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = "-"
$ws.Range("O14").Value = "-"

$ws.Range("P1").Select()
